# Fix InmemoryDB Primary Key's, Add SearchEngine and a couple of fixes in the header.
#
# The "Primary Key" column (B) of every little lookup table on the sheet was
# 0-based; bump every value in that column by 1 so the keys start at 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$primaryKeyFixes = [ordered]@{
    "B4"   = 1
    "B5"   = 2
    "B6"   = 3
    "B9"   = 1
    "B10"  = 2
    "B13"  = 1
    "B14"  = 2
    "B15"  = 3
    "B16"  = 4
    "B19"  = 1
    "B20"  = 2
    "B21"  = 3
    "B24"  = 1
    "B25"  = 2
    "B26"  = 3
    "B29"  = 1
    "B30"  = 2
    "B31"  = 3
    "B32"  = 4
    "B35"  = 1
    "B36"  = 2
    "B37"  = 3
    "B40"  = 1
    "B41"  = 2
    "B44"  = 1
    "B45"  = 2
    "B46"  = 3
    "B47"  = 4
    "B50"  = 1
    "B51"  = 2
    "B52"  = 3
    "B55"  = 1
    "B56"  = 2
    "B57"  = 3
    "B60"  = 1
    "B61"  = 2
    "B64"  = 1
    "B65"  = 2
    "B66"  = 3
    "B69"  = 1
    "B70"  = 2
    "B71"  = 3
    "B74"  = 1
    "B75"  = 2
    "B78"  = 1
    "B79"  = 2
    "B80"  = 3
    "B83"  = 1
    "B84"  = 2
    "B85"  = 3
    "B88"  = 1
    "B89"  = 2
    "B90"  = 3
    "B91"  = 4
    "B94"  = 1
    "B95"  = 2
    "B98"  = 1
    "B99"  = 2
    "B100" = 3
    "B103" = 1
    "B104" = 2
    "B105" = 3
    "B106" = 4
    "B109" = 1
    "B110" = 2
    "B113" = 1
    "B114" = 2
}

foreach ($addr in $primaryKeyFixes.Keys) {
    $ws.Range($addr).Value = $primaryKeyFixes[$addr]
}

# Match the author's final on-screen scroll position / selection (the sheet
# was left scrolled down near the bottom, with B118 selected).
$ws.Range("B118").Select()
